$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format temporarily so numeric-looking strings
# (e.g. "14.90", "36.191.77") are preserved exactly as text, not converted
# to floating point numbers by Excel's automatic type detection.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '36.191.77'
$ws.Range("E2").Value = '  -1.17%  '
$ws.Range("D3").Value = '2.014.84'
$ws.Range("E3").Value = '  -1.64%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").Value = '252.99'
$ws.Range("E5").Value = '  +3.05%  '
$ws.Range("D6").Value = '0.643'
$ws.Range("E6").Value = '  -3.67%  '
$ws.Range("D7").Value = '62.13'
$ws.Range("E7").Value = '  +9.12%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").Value = '59.06'
$ws.Range("E9").Value = '  -6.42%  '
$ws.Range("D10").Value = '0.371'
$ws.Range("E10").Value = '  +0.80%  '
$ws.Range("D11").Value = '0.0748'
$ws.Range("E11").Value = '  -0.53%  '
$ws.Range("E12").Value = '  -1.81%  '
$ws.Range("D13").Value = '0.918'
$ws.Range("E13").Value = '  -0.40%  '
$ws.Range("D14").Value = '14.90'
$ws.Range("E14").Value = '  +1.65%  '
$ws.Range("D15").Value = '2.312.78'
$ws.Range("E15").Value = '  -1.57%  '
$ws.Range("D16").Value = '5.43'
$ws.Range("E16").Value = '  -0.13%  '
$ws.Range("D17").Value = '19.55'
$ws.Range("E17").Value = '  +10.55%  '
$ws.Range("D18").Value = '2.015.06'
$ws.Range("E18").Value = '  -1.78%  '
$ws.Range("D19").Value = '36.130.35'
$ws.Range("E19").Value = '  -1.13%  '
$ws.Range("D20").Value = '72.11'
$ws.Range("E20").Value = '  +0.19%  '
$ws.Range("D21").Value = '0.0₃0859'
$ws.Range("E21").Value = '  -0.01%  '
$ws.Range("D22").Value = '5.28'
$ws.Range("E22").Value = '  +1.64%  '
$ws.Range("D23").Value = '234.26'
$ws.Range("E23").Value = '  -1.53%  '
$ws.Range("D24").Value = '2.70'
$ws.Range("E24").Value = '  +18.90%  '
$ws.Range("E25").Value = '  -0.09%  '
$ws.Range("D26").Value = '2.31'
$ws.Range("E26").Value = '  -2.28%  '
$ws.Range("D27").Value = '9.52'
$ws.Range("E27").Value = '  +2.43%  '
$ws.Range("D28").Value = '164.44'
$ws.Range("E28").Value = '  -0.25%  '
$ws.Range("D29").Value = '19.65'
$ws.Range("E29").Value = '  -1.79%  '
$ws.Range("D30").Value = '0.120'
$ws.Range("E30").Value = '  -0.99%  '
$ws.Range("D31").Value = '1.20'
$ws.Range("E31").Value = '  +0.14%  '
$ws.Range("E32").Value = '  +1.10%  '
$ws.Range("D33").Value = '0.109'
$ws.Range("E33").Value = '  +25.71%  '
$ws.Range("D34").Value = '0.0605'
$ws.Range("E34").Value = '  +0.25%  '
$ws.Range("D35").Value = '2.51'
$ws.Range("E35").Value = '  +13.24%  '
$ws.Range("D36").Value = '4.50'
$ws.Range("E36").Value = '  +1.37%  '
$ws.Range("E37").Value = '  -0.02%  '
$ws.Range("D38").Value = '1.81'
$ws.Range("E38").Value = '  -0.81%  '
$ws.Range("D39").Value = '5.83'
$ws.Range("E39").Value = '  +14.41%  '
$ws.Range("D40").Value = '0.102'
$ws.Range("E40").Value = '  +12.73%  '
$ws.Range("E41").Value = '  -0.89%  '
$ws.Range("E42").Value = '  -0.63%  '
$ws.Range("E43").Value = '  -0.27%  '
$ws.Range("E44").Value = '  +1.43%  '
$ws.Range("D45").Value = '16.73'
$ws.Range("E45").Value = '  +4.74%  '
$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").Value = '7.84'
$ws.Range("E46").Value = '  +4.73%  '
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = '93.92'
$ws.Range("E47").Value = '  -0.19%  '
$ws.Range("D48").Value = '1.424.27'
$ws.Range("E48").Value = '  +3.27%  '
$ws.Range("D49").Value = '2.51'
$ws.Range("E49").Value = '  +10.91%  '
$ws.Range("D50").Value = '2.90'
$ws.Range("E50").Value = '  -1.23%  '
$ws.Range("D51").Value = '47.57'
$ws.Range("E51").Value = '  +3.71%  '

# Restore the original (default) style on column D now that the values
# are safely stored as text, so no visible style/format change remains.
$ws.Range("D2:D51").Style = "Normal"

